# TradingModel - 2021/11/15 data update
# Refresh the open-position / today-close data table (rows 2-10) with the
# latest values. Stock_Id (col B) and TodayClose (col C) are refreshed for
# each position; col A holds the original source row index.
#
# Some TodayClose values must stay as literal text (matching how the sheet
# already stores a few prices, e.g. "267.00"/"165.50"/"264.00", as text
# rather than numbers). We force those with a leading apostrophe (classic
# "store as text" entry) and then strip the resulting quote-prefix
# formatting so the cell keeps its original (default) style - only its
# value/type changes, same as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $ws.Range($cell).Value = "'" + $text
    $ws.Range($cell).ClearFormats()
}

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2436
$ws.Range("C2").Value = 108

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 3035
$ws.Range("C3").Value = 198.5

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 3122
Set-TextValue "C4" "66.80"

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 3141
Set-TextValue "C5" "251.50"

# Row 6
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = 3588
$ws.Range("C6").Value = 163

# Row 7
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = 6104
Set-TextValue "C7" "182.00"

# Row 8
$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 6138
Set-TextValue "C8" "213.00"

# Row 9
$ws.Range("A9").Value = 13
$ws.Range("B9").Value = 6271
$ws.Range("C9").Value = 302.5

# Row 10
$ws.Range("A10").Value = 14
$ws.Range("B10").Value = 6411
Set-TextValue "C10" "290.00"
